{"js": "// Blog-content fixes for \"week 15\" doc, per commit \":sparkle: fixed blog contents\".\n//\n// Four textual edits are applied (paragraph-scoped search & replace so we never\n// touch a similar phrase living in a different paragraph):\n//   1. \"...We had a fun Kahoot! game that made me...\"\n//        -> \"...We had a fun Kahoot! A game that made me...\"\n//   2. \"...that my knowledge of referencing and citations was limited.\" (unchanged\n//      visible text - only the run layout changed in the source - nothing to do)\n//   3. \"...that we understand the extent of the published work...\"\n//        -> \"...that we fully understand the published work...\"\n//   4. \"...exploring concepts on recycling ... I also started working on exploring\n//        waste classifier models online and using an ESP32...\"\n//        -> \"...explored concepts on recycling ... I explored waste classifier\n//        models online and used an ESP32...\"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfunction findParagraph(substring) {\n  const p = paragraphs.items.find((item) => item.text.indexOf(substring) !== -1);\n  if (!p) {\n    throw new Error(\"Could not locate paragraph containing: \" + substring);\n  }\n  return p;\n}\n\nasync function replaceOnce(scopeRange, findText, replaceText) {\n  const results = scopeRange.search(findText, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Could not locate text to replace: \" + findText);\n  }\n  results.items[0].insertText(replaceText, \"Replace\");\n  await context.sync();\n}\n\n// --- Edit 1: \"Kahoot! game\" -> \"Kahoot! A game\" -----------------------------\nconst p1 = findParagraph(\"We had a fun Kahoot! game that made me\");\nawait replaceOnce(p1, \"Kahoot! game\", \"Kahoot! A game\");\n\n// --- Edit 3: \"understand the extent of the published work\" -> \"fully understand the published work\"\nconst p3 = findParagraph(\"we understand the extent of the published work\");\nawait replaceOnce(\n  p3,\n  \"we understand the extent of the published work\",\n  \"we fully understand the published work\"\n);\n\n// --- Edit 4: the ESP32 / literature review reflection paragraph ------------\nconst p4 = findParagraph(\"using an ESP32 WIFI + Camera module\");\nawait replaceOnce(p4, \"exploring concepts on recycling\", \"explored concepts on recycling\");\nawait replaceOnce(\n  p4,\n  \"I also started working on exploring waste classifier models online and using an ESP32\",\n  \"I explored waste classifier models online and used an ESP32\"\n);\n", "ps1": "# Blog-content fixes for \"week 15\" doc, per commit \":sparkle: fixed blog contents\".\n#\n# Four textual edits are applied as plain-text Find/Replace passes, each scoped to\n# the paragraph containing the target phrase (so we never touch a similar phrase\n# living in a different paragraph):\n#   1. \"...We had a fun Kahoot! game that made me...\"\n#        -> \"...We had a fun Kahoot! A game that made me...\"\n#   2. \"...that my knowledge of referencing and citations was limited.\" (unchanged\n#      visible text - only the run layout changed in the source - nothing to do)\n#   3. \"...that we understand the extent of the published work...\"\n#        -> \"...that we fully understand the published work...\"\n#   4. \"...exploring concepts on recycling ... I also started working on exploring\n#        waste classifier models online and using an ESP32...\"\n#        -> \"...explored concepts on recycling ... I explored waste classifier\n#        models online and used an ESP32...\"\n\n$d = $word.ActiveDocument\n\nfunction Get-ParagraphContaining($substring) {\n    foreach ($p in $d.Paragraphs) {\n        if ($p.Range.Text -like \"*$substring*\") {\n            return $p\n        }\n    }\n    throw \"Could not locate paragraph containing: $substring\"\n}\n\nfunction Replace-InParagraph($paragraph, $findText, $replaceText) {\n    $rng = $paragraph.Range\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $ok = $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $ok) {\n        throw \"Could not locate text to replace: $findText\"\n    }\n}\n\n# --- Edit 1: \"Kahoot! game\" -> \"Kahoot! A game\" -----------------------------\n$p1 = Get-ParagraphContaining \"We had a fun Kahoot! game that made me\"\nReplace-InParagraph $p1 \"Kahoot! game\" \"Kahoot! A game\"\n\n# --- Edit 3: \"understand the extent of the published work\" -> \"fully understand the published work\"\n$p3 = Get-ParagraphContaining \"we understand the extent of the published work\"\nReplace-InParagraph $p3 \"we understand the extent of the published work\" \"we fully understand the published work\"\n\n# --- Edit 4: the ESP32 / literature review reflection paragraph ------------\n$p4 = Get-ParagraphContaining \"using an ESP32 WIFI + Camera module\"\nReplace-InParagraph $p4 \"exploring concepts on recycling\" \"explored concepts on recycling\"\nReplace-InParagraph $p4 \"I also started working on exploring waste classifier models online and using an ESP32\" \"I explored waste classifier models online and used an ESP32\"\n"}
